$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.223.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.489.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.39%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.76%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").Value = "  -0.60%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.490.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.137"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.64%  "

$ws.Range("E11").Value = "  +0.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.38%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.333"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.940.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.15%  "

$ws.Range("E15").Value = "  -1.64%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.140.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.47%  "

$ws.Range("E17").Value = "  -1.64%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.479.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.72%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "349.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.43%  "

$ws.Range("E23").Value = "  -0.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.19%  "

$ws.Range("E25").Value = "  -3.97%  "

$ws.Range("E26").Value = "  -2.31%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.23%  "

$ws.Range("E28").Value = "  -0.12%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.619.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.60%  "

$ws.Range("E30").Value = "  -3.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "509.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.45%  "

$ws.Range("E33").Value = "  -2.89%  "

$ws.Range("E34").Value = "  -3.76%  "

$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.48%  "

$ws.Range("E37").Value = "  -7.17%  "

$ws.Range("E38").Value = "  +0.89%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.24"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.91%  "

$ws.Range("E40").Value = "  -5.37%  "

$ws.Range("E41").Value = "  -2.32%  "

$ws.Range("E42").Value = "  -0.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.328"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.86%  "

$ws.Range("E45").Value = "  -3.70%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.85"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.72%  "

$ws.Range("E48").Value = "  -3.91%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.514"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.41%  "

$ws.Range("E50").Value = "  -6.07%  "

$ws.Range("E51").Value = "  -0.71%  "
